$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-25: Url, Passed, Comments reflect the new availability check
# (A "Key", C "Page", D "Test Case" are unchanged)

$ws.Range("B2").Value = 'https://www.varoom.com/property/anurra1-domus-de-janas/BC-4871160'
$ws.Range("E2").Value = $False
$ws.Range("F2").Value = 'The property ''Anurra1-Domus de Janas'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B3").Value = 'https://www.varoom.com/property/nice-apartment-for-5-guests-with-tv-terrace-and-pets-allowed/EP-27689065'
$ws.Range("E3").Value = $False
$ws.Range("F3").Value = 'The property ''Nice apartment for 5 guests with TV, terrace and pets allowed'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B4").Value = 'https://www.varoom.com/property/holiday-home-isola-rossa/BC-6119796'
$ws.Range("E4").Value = $False
$ws.Range("F4").Value = 'The property ''Holiday Home Isola Rossa'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B5").Value = 'https://www.varoom.com/property/costa-paradiso-resort/EP-102231382'
$ws.Range("E5").Value = $False
$ws.Range("F5").Value = 'The property ''Costa Paradiso Resort'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B6").Value = 'https://www.varoom.com/property/locazione-turistica-scalitti-by-interhome/BC-3551349'
$ws.Range("E6").Value = $False
$ws.Range("F6").Value = 'The property ''Locazione Turistica Scalitti by Interhome'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B7").Value = 'https://www.varoom.com/property/isola-rossa-borgo-mare-agenzia-isola-rossa/EP-11703078'
$ws.Range("E7").Value = $False
$ws.Range("F7").Value = 'The property ''Isola Rossa Borgo Mare - Agenzia Isola Rossa'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B8").Value = 'https://www.varoom.com/property/walk-to-the-beach-from-your-cottage-apartment-set-in-wild-rural-sardinia/BC-3887684'
$ws.Range("E8").Value = $False
$ws.Range("F8").Value = 'The property ''Walk To The Beach From Your Cottage-Apartment Set In Wild, Rural Sardinia'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B9").Value = 'https://www.varoom.com/property/gravina-resort-apartments/EP-5175840'
$ws.Range("E9").Value = $False
$ws.Range("F9").Value = 'The property ''Gravina Resort & Apartments'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B10").Value = 'https://www.varoom.com/property/affittimoderni-isola-rossa-borgo/BC-4824656'
$ws.Range("E10").Value = $False
$ws.Range("F10").Value = 'The property ''Affittimoderni Isola Rossa Borgo'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B11").Value = 'https://www.varoom.com/property/residence-with-pool-in-isola-rossa/EP-100779242'
$ws.Range("E11").Value = $False
$ws.Range("F11").Value = 'The property ''Residence with pool in Isola Rossa'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B12").Value = 'https://www.varoom.com/property/cottage-apartment-in-rural-sardinia-with-sun-sea-and-sand/BC-3765161'
$ws.Range("E12").Value = $False
$ws.Range("F12").Value = 'The property ''Cottage-Apartment In Rural Sardinia With Sun, Sea And Sand'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B13").Value = 'https://www.varoom.com/property/villa-rocce-rosse-costa-paradiso/EP-93506023'
$ws.Range("E13").Value = $False
$ws.Range("F13").Value = 'The property ''Villa Rocce Rosse Costa Paradiso'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B14").Value = 'https://www.varoom.com/property/il-boschetto-dei-corbezzoli-villetta-3/BC-5288835'
$ws.Range("E14").Value = $False
$ws.Range("F14").Value = 'The property ''Il Boschetto dei Corbezzoli villetta 3'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B15").Value = 'https://www.varoom.com/property/apartment-with-stunning-views/EP-96818551'
$ws.Range("E15").Value = $False
$ws.Range("F15").Value = 'The property ''Apartment With Stunning Views'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B16").Value = 'https://www.varoom.com/property/costa-paradiso-villaggio-tamerici/BC-1065581'
$ws.Range("E16").Value = $False
$ws.Range("F16").Value = 'The property ''Costa Paradiso Villaggio Tamerici'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B17").Value = 'https://www.varoom.com/property/villetta-dei-ginepri-costa-paradiso/EP-94632555'
$ws.Range("E17").Value = $False
$ws.Range("F17").Value = 'The property ''Villetta dei Ginepri Costa Paradiso'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B18").Value = 'https://www.varoom.com/property/il-boschetto-dei-corbezzoli-villetta-2/BC-5283284'
$ws.Range("E18").Value = $False
$ws.Range("F18").Value = 'The property ''Il Boschetto dei Corbezzoli Villetta 2'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B19").Value = 'https://www.varoom.com/property/isola-rossa-apartment-with-breathtaking-sea-view/EP-96814262'
$ws.Range("E19").Value = $False
$ws.Range("F19").Value = 'The property ''Isola Rossa Apartment With Breathtaking sea View'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B20").Value = 'https://www.varoom.com/property/scoglio-a-costa-paradiso/BC-3898644'
$ws.Range("E20").Value = $False
$ws.Range("F20").Value = 'The property ''Scoglio a Costa Paradiso'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B21").Value = 'https://www.varoom.com/property/fronte-mare/BC-2492992'
$ws.Range("E21").Value = $False
$ws.Range("F21").Value = 'The property ''Fronte Mare'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B22").Value = 'https://www.varoom.com/property/apartment-with-swimming-pool-in-trinit-d-agultu-e-vignola/BC-2469706'
$ws.Range("E22").Value = $False
$ws.Range("F22").Value = 'The property ''Apartment with Swimming Pool in Trinit d Agultu e Vignola'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B23").Value = 'https://www.varoom.com/property/central-apartment-irina-with-terrace/BC-2192182'
$ws.Range("E23").Value = $False
$ws.Range("F23").Value = 'The property ''Central apartment Irina with terrace'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B24").Value = 'https://www.varoom.com/property/camera-con-bagno-isola-rossa-paduledda/BC-8957107'
$ws.Range("E24").Value = $False
$ws.Range("F24").Value = 'The property ''Camera con bagno Isola Rossa, Paduledda'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("B25").Value = 'https://www.varoom.com/property/appartamento-isola-rossa-paduledda/BC-8879596'
$ws.Range("E25").Value = $False
$ws.Range("F25").Value = 'The property ''Appartamento Isola Rossa, Paduledda'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

# Append new rows 26-37 with the same result pattern

$ws.Range("A26").Value = 'www.varoom.com'
$ws.Range("B26").Value = 'https://www.varoom.com/property/sweet-costa-paradiso-splendida-vista-mare/BC-8834419'
$ws.Range("C26").Value = 'Hybrid'
$ws.Range("D26").Value = 'Property available in date range'
$ws.Range("E26").Value = $False
$ws.Range("F26").Value = 'The property ''Sweet Costa Paradiso splendida vista mare'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("A27").Value = 'www.varoom.com'
$ws.Range("B27").Value = 'https://www.varoom.com/property/appartamento-vista-mare-via-tinnari/BC-8849180'
$ws.Range("C27").Value = 'Hybrid'
$ws.Range("D27").Value = 'Property available in date range'
$ws.Range("E27").Value = $False
$ws.Range("F27").Value = 'The property ''Appartamento Vista Mare via Tinnari'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("A28").Value = 'www.varoom.com'
$ws.Range("B28").Value = 'https://www.varoom.com/property/emanuele-villetta-con-ampio-giardino-piscina-condizionatori-caldo-freddo/BC-7890193'
$ws.Range("C28").Value = 'Hybrid'
$ws.Range("D28").Value = 'Property available in date range'
$ws.Range("E28").Value = $False
$ws.Range("F28").Value = 'The property ''Emanuele Villetta con Ampio Giardino Piscina CONDIZIONATORI CALDO FREDDO'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("A29").Value = 'www.varoom.com'
$ws.Range("B29").Value = 'https://www.varoom.com/property/costa-paradiso-comprensorio-in-totale-relax/BC-7635836'
$ws.Range("C29").Value = 'Hybrid'
$ws.Range("D29").Value = 'Property available in date range'
$ws.Range("E29").Value = $False
$ws.Range("F29").Value = 'The property ''Costa paradiso comprensorio in totale relax'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("A30").Value = 'www.varoom.com'
$ws.Range("B30").Value = 'https://www.varoom.com/property/dream-isola-rossa/BC-7220073'
$ws.Range("C30").Value = 'Hybrid'
$ws.Range("D30").Value = 'Property available in date range'
$ws.Range("E30").Value = $False
$ws.Range("F30").Value = 'The property ''Dream Isola Rossa'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("A31").Value = 'www.varoom.com'
$ws.Range("B31").Value = 'https://www.varoom.com/property/stazzo-jana/BC-11063229'
$ws.Range("C31").Value = 'Hybrid'
$ws.Range("D31").Value = 'Property available in date range'
$ws.Range("E31").Value = $False
$ws.Range("F31").Value = 'The property ''Stazzo Jana'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("A32").Value = 'www.varoom.com'
$ws.Range("B32").Value = 'https://www.varoom.com/property/anurra1-domus-de-janas/BC-4871160'
$ws.Range("C32").Value = 'Hybrid'
$ws.Range("D32").Value = 'Property available in date range'
$ws.Range("E32").Value = $False
$ws.Range("F32").Value = 'The property ''Anurra1-Domus de Janas'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("A33").Value = 'www.varoom.com'
$ws.Range("B33").Value = 'https://www.varoom.com/property/nice-apartment-for-5-guests-with-tv-terrace-and-pets-allowed/EP-27689065'
$ws.Range("C33").Value = 'Hybrid'
$ws.Range("D33").Value = 'Property available in date range'
$ws.Range("E33").Value = $False
$ws.Range("F33").Value = 'The property ''Nice apartment for 5 guests with TV, terrace and pets allowed'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("A34").Value = 'www.varoom.com'
$ws.Range("B34").Value = 'https://www.varoom.com/property/holiday-home-isola-rossa/BC-6119796'
$ws.Range("C34").Value = 'Hybrid'
$ws.Range("D34").Value = 'Property available in date range'
$ws.Range("E34").Value = $False
$ws.Range("F34").Value = 'The property ''Holiday Home Isola Rossa'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("A35").Value = 'www.varoom.com'
$ws.Range("B35").Value = 'https://www.varoom.com/property/costa-paradiso-resort/EP-102231382'
$ws.Range("C35").Value = 'Hybrid'
$ws.Range("D35").Value = 'Property available in date range'
$ws.Range("E35").Value = $False
$ws.Range("F35").Value = 'The property ''Costa Paradiso Resort'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("A36").Value = 'www.varoom.com'
$ws.Range("B36").Value = 'https://www.varoom.com/property/locazione-turistica-scalitti-by-interhome/BC-3551349'
$ws.Range("C36").Value = 'Hybrid'
$ws.Range("D36").Value = 'Property available in date range'
$ws.Range("E36").Value = $False
$ws.Range("F36").Value = 'The property ''Locazione Turistica Scalitti by Interhome'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'

$ws.Range("A37").Value = 'www.varoom.com'
$ws.Range("B37").Value = 'https://www.varoom.com/property/isola-rossa-borgo-mare-agenzia-isola-rossa/EP-11703078'
$ws.Range("C37").Value = 'Hybrid'
$ws.Range("D37").Value = 'Property available in date range'
$ws.Range("E37").Value = $False
$ws.Range("F37").Value = 'The property ''Isola Rossa Borgo Mare - Agenzia Isola Rossa'' is Unavailable in the specified date range. | Location: Shawnview, Dates: 2025-02-01 to 2025-02-04'
